$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1. Fill in the sample sizes for the "Cushing's (...)" sentence.
#    This replace spans the old, now-stale grammar-check markers
#    (proofErr gramStart/gramEnd) so they get swept away along with
#    the text they used to bracket.
# ------------------------------------------------------------------
$d.Content.Find.Execute(
    "non-obese n=; obese n= ) and",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "non-obese n=3; obese n=5 ) and", 2) | Out-Null

# ------------------------------------------------------------------
# 2. Fill in the sample sizes for the "control (...)" sentence.
# ------------------------------------------------------------------
$d.Content.Find.Execute(
    "non-obese n=; obese n=) BMI",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "non-obese n=5; obese n=6) BMI", 2) | Out-Null

# ------------------------------------------------------------------
# 3. Remove the old review comment (id=0) that asked for the sample
#    sizes filled in above -- it is now resolved.  Comment.Delete()
#    also removes the commentRangeStart/commentRangeEnd markers and
#    the commentReference run bound to it.
# ------------------------------------------------------------------
$d.Comments.Item(1).Delete()

# ------------------------------------------------------------------
# 4. Move the "_GoBack" bookmark: drop it from its old spot near the
#    end of the document and re-create it right after the sample-size
#    text that was just typed in (where the resolved comment used to
#    close).
# ------------------------------------------------------------------
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}

$anchor = $d.Content
$anchor.Find.Execute(
    "non-obese n=5; obese n=6)",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "", 0) | Out-Null
$bmRange = $anchor.Duplicate
$bmRange.Collapse(0)
$d.Bookmarks.Add("_GoBack", $bmRange) | Out-Null
